$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New video IDs to be appended after deleting the stale entry (row 18)
# and shifting everything up, simulating deleter.py removing a watched
# video and converter.py appending freshly converted ones.

# Delete the stale video row (old A18 = "lFlWULCpmak"); this shifts the
# rows below it up by one.
$ws.Range("A18").EntireRow.Delete()

# Append the newly converted video IDs at the bottom of the list.
$newVideoIds = @(
    "dhZUsNJ-LQU",
    "j2r2nDhTzO4",
    "oadhHk2xs6c",
    "UPOVM_oYxHc",
    "y_KCK-pHzqk",
    "eUKhgjTtxyM",
    "OzWrVeC-GGw",
    "IoCcF0UrQOQ",
    "G5weJd_FwAo",
    "QCUwbIQIP8E",
    "0MIXDyQAjAE",
    "jWeFH9QyLRE",
    "kujwJhXRGSs",
    "VYSc1h8qkgg",
    "z-DySQ5PAAc",
    "mPcLc9qgBS8",
    "2cTXeSVrSD0",
    "iOsE0eANCmA",
    "gCWj8Nz5DUg",
    "6WGB6lK6pAc"
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
for ($i = 0; $i -lt $newVideoIds.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newVideoIds[$i]
}
